$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Address URL cell to point to the new /rest/ endpoint
$ws.Range("A5").Value = "http://schattenhauer.de/GiveMeTheRESTServer/rest/ FUNCTION"

# Move the active selection from C4 to A5
$ws.Range("A5").Select()
